{"js": "// Replace the date heading and the twenty-five \"NNN\u00d7N=\" problems in the\n// multiplication worksheet with their updated values (see commit diff).\nconst replacements = [\n  [\"2025-05-15 Thursday\", \"2025-05-16 Friday\"],\n  [\"767\u00d78=\", \"788\u00d73=\"],\n  [\"590\u00d74=\", \"357\u00d76=\"],\n  [\"698\u00d73=\", \"343\u00d72=\"],\n  [\"181\u00d72=\", \"935\u00d76=\"],\n  [\"233\u00d76=\", \"365\u00d72=\"],\n  [\"482\u00d77=\", \"393\u00d74=\"],\n  [\"832\u00d74=\", \"784\u00d76=\"],\n  [\"134\u00d73=\", \"296\u00d72=\"],\n  [\"176\u00d77=\", \"692\u00d78=\"],\n  [\"215\u00d72=\", \"371\u00d73=\"],\n  [\"452\u00d74=\", \"260\u00d75=\"],\n  [\"403\u00d74=\", \"979\u00d75=\"],\n  [\"728\u00d78=\", \"954\u00d78=\"],\n  [\"856\u00d79=\", \"824\u00d75=\"],\n  [\"892\u00d73=\", \"863\u00d77=\"],\n  [\"715\u00d73=\", \"718\u00d72=\"],\n  [\"225\u00d76=\", \"256\u00d75=\"],\n  [\"489\u00d76=\", \"831\u00d76=\"],\n  [\"832\u00d77=\", \"299\u00d79=\"],\n  [\"260\u00d76=\", \"371\u00d74=\"],\n  [\"739\u00d72=\", \"171\u00d78=\"],\n  [\"287\u00d79=\", \"122\u00d73=\"],\n  [\"261\u00d74=\", \"190\u00d76=\"],\n  [\"475\u00d72=\", \"643\u00d78=\"],\n  [\"469\u00d76=\", \"336\u00d73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date heading and the twenty-five \"NNN\u00d7N=\" problems in the\n# multiplication worksheet with their updated values (see commit diff).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-05-15 Thursday\", \"2025-05-16 Friday\"),\n    @(\"767\u00d78=\", \"788\u00d73=\"),\n    @(\"590\u00d74=\", \"357\u00d76=\"),\n    @(\"698\u00d73=\", \"343\u00d72=\"),\n    @(\"181\u00d72=\", \"935\u00d76=\"),\n    @(\"233\u00d76=\", \"365\u00d72=\"),\n    @(\"482\u00d77=\", \"393\u00d74=\"),\n    @(\"832\u00d74=\", \"784\u00d76=\"),\n    @(\"134\u00d73=\", \"296\u00d72=\"),\n    @(\"176\u00d77=\", \"692\u00d78=\"),\n    @(\"215\u00d72=\", \"371\u00d73=\"),\n    @(\"452\u00d74=\", \"260\u00d75=\"),\n    @(\"403\u00d74=\", \"979\u00d75=\"),\n    @(\"728\u00d78=\", \"954\u00d78=\"),\n    @(\"856\u00d79=\", \"824\u00d75=\"),\n    @(\"892\u00d73=\", \"863\u00d77=\"),\n    @(\"715\u00d73=\", \"718\u00d72=\"),\n    @(\"225\u00d76=\", \"256\u00d75=\"),\n    @(\"489\u00d76=\", \"831\u00d76=\"),\n    @(\"832\u00d77=\", \"299\u00d79=\"),\n    @(\"260\u00d76=\", \"371\u00d74=\"),\n    @(\"739\u00d72=\", \"171\u00d78=\"),\n    @(\"287\u00d79=\", \"122\u00d73=\"),\n    @(\"261\u00d74=\", \"190\u00d76=\"),\n    @(\"475\u00d72=\", \"643\u00d78=\"),\n    @(\"469\u00d76=\", \"336\u00d73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
